# Applies odds updates to Sheet1, matching the commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 - Welsh Premiership: Connahs Quay vs The New Saints
$ws.Range("G6").Value = 7.8
$ws.Range("H6").Value = 1.45
$ws.Range("K6").Value = 6.4

# Row 7 - Saudi Professional League: NEOM Sports Club vs Al-Ittihad
$ws.Range("H7").Value = 2.02
$ws.Range("P7").Value = 2.64

# Row 8 - Israeli Premier League: Maccabi Netanya vs Hapoel Petach Tikva
$ws.Range("G8").Value = 2.34
$ws.Range("J8").Value = 3.35
$ws.Range("W8").Value = 1.75

# Row 9 - Saudi Professional League: Al-Shabab (KSA) vs Al-Quadisiya (KSA)
$ws.Range("AK9").Value = 90
$ws.Range("AL9").Value = 80

# Row 10 - Saudi Professional League: Al-Kholood Club vs Al-Hilal
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 1.25
$ws.Range("J10").Value = 6.6
$ws.Range("K10").Value = 8
$ws.Range("Q10").Value = 1.25
$ws.Range("T10").Value = 1.04
$ws.Range("W10").Value = 1.07
